$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 'Datos actualizados a 26 de Mayo de 2020 a las 00:05'

$ws.Range("B4").Value = 1703881
$ws.Range("C4").Value = 17445
$ws.Range("D4").Value = 461002
$ws.Range("E4").Value = 1143125
$ws.Range("G4").Value = 454
$ws.Range("H4").Value = 99754

$ws.Range("D31").Value = 28200
$ws.Range("E31").Value = 633

$ws.Range("A36").Value = 'Colombia'
$ws.Range("B36").Value = 21981
$ws.Range("C36").Value = 806
$ws.Range("D36").Value = 5265
$ws.Range("E36").Value = 15966
$ws.Range("G36").Value = 23
$ws.Range("H36").Value = 750

$ws.Range("A37").Value = 'Kuwait'
$ws.Range("B37").Value = 21967
$ws.Range("C37").Value = 665
$ws.Range("D37").Value = 6621
$ws.Range("E37").Value = 15181
$ws.Range("G37").Value = 9
$ws.Range("H37").Value = 165

$ws.Range("A38").Value = 'Polonia'
$ws.Range("B38").Value = 21631
$ws.Range("C38").Value = 305
$ws.Range("D38").Value = 9276
$ws.Range("E38").Value = 11348
$ws.Range("G38").Value = 11
$ws.Range("H38").Value = 1007

$ws.Range("A39").Value = 'Ucrania'
$ws.Range("B39").Value = 21245
$ws.Range("C39").Value = 259
$ws.Range("D39").Value = 7234
$ws.Range("E39").Value = 13388
$ws.Range("G39").Value = 6
$ws.Range("H39").Value = 623

$ws.Range("B141").Value = 386
$ws.Range("C141").Value = 5
$ws.Range("D141").Value = 161
$ws.Range("E141").Value = 212
$ws.Range("G141").Value = 1
$ws.Range("H141").Value = 13

$ws.Range("B154").Value = 222
$ws.Range("C154").Value = 24
$ws.Range("D154").Value = 69
$ws.Range("E154").Value = 153

$ws.Range("D166").Value = 89
$ws.Range("E166").Value = 35

$ws.Range("A197").Value = 'Curazao'
$ws.Range("C197").Value = 1
$ws.Range("D197").Value = 14
$ws.Range("H197").Value = 1

$ws.Range("A198").Value = 'Fiyi'
$ws.Range("C198").Value = 0
$ws.Range("D198").Value = 15
$ws.Range("H198").Value = 0

$ws.Range("A199").Value = 'Nueva Caledonia'

$ws.Range("A201").Value = 'Santa Lucia'

$ws.Range("A215").Value = 'Bonaire, San Eustaquio y Saba'

$ws.Range("A216").Value = 'San Bartolome'
